$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Row 1: rename test case label to the new positive-test variant
$ws.Range("A1").Value = "AutoInsuranceStartQuoteByZipcode_PositiveTest"

# Update selection to E7 (was A7)
$ws.Range("E7").Select()

# Widen column A
$ws.Columns.Item(1).ColumnWidth = 47.6

# Row 8: turn it into a duplicate header row (same pattern as row 2), adding I8/J8
$ws.Range("E8").Value = $ws.Range("E2").Value2
$ws.Range("F8").Value = $ws.Range("F2").Value2
$ws.Range("G8").Value = $ws.Range("G2").Value2
$ws.Range("H8").Value = $ws.Range("H2").Value2
$ws.Range("I8").Value = $ws.Range("I2").Value2
$ws.Range("J8").Value = $ws.Range("J2").Value2

# Rows 9-12: drop the sample data in B,C,D,F,G and blank out the date cells E,H
# (keeping their existing date number-format style)
$ws.Range("B9:D9").ClearContents()
$ws.Range("F9:G9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("H9").ClearContents()

$ws.Range("B10:D10").ClearContents()
$ws.Range("F10:G10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("H10").ClearContents()

$ws.Range("B11:D11").ClearContents()
$ws.Range("F11:G11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("H11").ClearContents()

$ws.Range("B12:D12").ClearContents()
$ws.Range("F12:G12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("H12").ClearContents()

$wb.Save()
